$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.983274459838867
$ws.Range("B1").Value = 4.185163974761963
$ws.Range("C1").Value = 2.933164834976196
$ws.Range("D1").Value = 2.332669496536255
$ws.Range("E1").Value = 1.936059236526489
